$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Spa Endurance" row data (G4/H4 lookup pair + new row 6 totals)
$ws.Range("G4").Value = "Spa Endurance"
$ws.Range("H4").Value = 7.004

$ws.Range("B6").Value = 553
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = 123
$ws.Range("E6").Formula = "=B6*H4"

$ws.Range("J7").Select()
